# fix upload rapel manfaat
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections on the transaction rows (column H = "Kode Pensiun") ---
# Row 11
$ws.Range("H11").Value = 435643

# Row 12
$ws.Range("H12").Value = 2

# Row 13: Kode Pensiun value fixed, and "Nama Peserta" (column I) corrected
$ws.Range("H13").Value = 4341
$ws.Range("I13").Value = "rfsg"

# --- Column C width: widen / best-fit to show full "Kode Voucher" text ---
$ws.Columns("C:C").ColumnWidth = 28.3

# --- Selection moved to F15 with the view scrolled so column B is leftmost ---
$ws.Range("F15").Select()
